$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old header row (row 2: "municipio deresidencia","Casos","Obitos")
# so that the data rows shift up by one.
$ws.Rows.Item(2).Delete()

# Remove the trailing "outros estados" / "outros paises" rows (now rows 50 and 51
# after the shift above).
$ws.Rows.Item(51).Delete()
$ws.Rows.Item(50).Delete()
